$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate new translation rows 415-462 (lab.vape.* form strings)
$ws.Range("A415").Value = 'cs'
$ws.Range("B415").Value = 'lab.vape.setupId.label'
$ws.Range("C415").Value = 'Setup'
$ws.Range("A416").Value = 'cs'
$ws.Range("B416").Value = 'lab.vape.create.submit'
$ws.Range("C416").Value = 'Vytvořit vape'
$ws.Range("A417").Value = 'cs'
$ws.Range("B417").Value = 'lab.setup.tooltip.create'
$ws.Range("C417").Value = 'Vytvořit setup'
$ws.Range("A418").Value = 'cs'
$ws.Range("B418").Value = 'lab.vape.mixtureId.label'
$ws.Range("C418").Value = 'Mix'
$ws.Range("A419").Value = 'cs'
$ws.Range("B419").Value = 'lab.mixture.tooltip.create'
$ws.Range("C419").Value = 'Vytvořit mix'
$ws.Range("A420").Value = 'cs'
$ws.Range("B420").Value = 'lab.mixture.create.success'
$ws.Range("C420").Value = 'Mix [{{data.name}}] byl úspěšně uložen.'
$ws.Range("A421").Value = 'cs'
$ws.Range("B421").Value = 'lab.setup.create.success'
$ws.Range("C421").Value = 'Setup [{{data.name}}] byl úspěšně uložen.'
$ws.Range("A422").Value = 'cs'
$ws.Range("B422").Value = 'lab.vape.driptipId.label'
$ws.Range("C422").Value = 'Náústek'
$ws.Range("A423").Value = 'cs'
$ws.Range("B423").Value = 'lab.driptip.name.label'
$ws.Range("C423").Value = 'Jméno'
$ws.Range("A424").Value = 'cs'
$ws.Range("B424").Value = 'lab.driptip.create.success'
$ws.Range("C424").Value = 'Náústek [{{data.name}}] byl uložen.'
$ws.Range("A425").Value = 'cs'
$ws.Range("B425").Value = 'lab.vape.common.title'
$ws.Range("C425").Value = 'Obecné'
$ws.Range("A426").Value = 'cs'
$ws.Range("B426").Value = 'lab.vape.rating.title'
$ws.Range("C426").Value = 'Celkové hodnocení'
$ws.Range("A427").Value = 'cs'
$ws.Range("B427").Value = 'lab.vape.rating.label'
$ws.Range("C427").Value = 'Celkové hodnocení'
$ws.Range("A428").Value = 'cs'
$ws.Range("B428").Value = 'lab.vape.rating.label.tooltip'
$ws.Range("C428").Value = 'Tato hodnota by měla reprezentovat celkový pocit z vapingu, včetně setupu, bublání atomizéru, chuti, prostě všeho.'
$ws.Range("A429").Value = 'cs'
$ws.Range("B429").Value = 'lab.vape.taste.label'
$ws.Range("C429").Value = 'Chuťový projev'
$ws.Range("A430").Value = 'cs'
$ws.Range("B430").Value = 'lab.vape.taste.label.tooltip'
$ws.Range("C430").Value = 'Tato hodnota sleduje pouze dojem z chuťového projevu buildu. Hlavní tedy je, jak moc je projev blízko očekávání.'
$ws.Range("A431").Value = 'cs'
$ws.Range("B431").Value = 'lab.vape.rating-advanced.title'
$ws.Range("C431").Value = 'Rozborka chuťového projevu'
$ws.Range("A432").Value = 'cs'
$ws.Range("B432").Value = 'lab.vape.fruits.label'
$ws.Range("C432").Value = 'Ovocné tóny'
$ws.Range("A433").Value = 'cs'
$ws.Range("B433").Value = 'lab.vape.fruits.label.tooltip'
$ws.Range("C433").Value = 'Tato položka sleduje, jak dobře vybraný build podává ovocné složky; např. v daném nastavení se mohou lépe projevovat citronové tóny, ale chuť buchty může být v pozadí.'
$ws.Range("A434").Value = 'cs'
$ws.Range("B434").Value = 'lab.vape.tobacco.label'
$ws.Range("C434").Value = 'Tabák'
$ws.Range("A435").Value = 'cs'
$ws.Range("B435").Value = 'lab.vape.tobacco.label.tooltip'
$ws.Range("C435").Value = 'Hodnocení tabákového projevu daného liquidu.'
$ws.Range("A436").Value = 'cs'
$ws.Range("B436").Value = 'lab.vape.cakes.label'
$ws.Range("C436").Value = 'Buchty'
$ws.Range("A437").Value = 'cs'
$ws.Range("B437").Value = 'lab.vape.cakes.label.tooltip'
$ws.Range("C437").Value = 'Hodnocení kvality podání chuti buchet.'
$ws.Range("A438").Value = 'cs'
$ws.Range("B438").Value = 'lab.vape.complex.label'
$ws.Range("C438").Value = 'Komplexní'
$ws.Range("A439").Value = 'cs'
$ws.Range("B439").Value = 'lab.vape.complex.label.tooltip'
$ws.Range("C439").Value = 'Toto hodnocení celkově sleduje komplexitu podání chuti u liquidů, kde není dominantní položka (např. pouze jablko); lze tak také určit, jak dobře daný build podává komplikovanější příchutě (např. tabák s tóny bourbonu).'
$ws.Range("A440").Value = 'cs'
$ws.Range("B440").Value = 'lab.vape.fresh.label'
$ws.Range("C440").Value = 'Větrnost'
$ws.Range("A441").Value = 'cs'
$ws.Range("B441").Value = 'lab.vape.fresh.label.tooltip'
$ws.Range("C441").Value = 'Liquidy, které obsahují mátu nebo jsou jinak větravé, nejsou vždy příjemné; toto hodnocení je reverzní - vyšší číslo udává větší míru "ice" efektu, kdy nejvyšší reprezentuje již nepříjemný zážitek.'
$ws.Range("A442").Value = 'cs'
$ws.Range("B442").Value = 'lab.vape.vape.title'
$ws.Range("C442").Value = 'Hodnocení vapingu'
$ws.Range("A443").Value = 'cs'
$ws.Range("B443").Value = 'lab.vape.settings.title'
$ws.Range("C443").Value = 'Nastavení'
$ws.Range("A444").Value = 'cs'
$ws.Range("B444").Value = 'lab.vape.power.label'
$ws.Range("C444").Value = 'Použitý výkon'
$ws.Range("A445").Value = 'cs'
$ws.Range("B445").Value = 'lab.vape.power.label.tooltip'
$ws.Range("C445").Value = 'Zde si zaznamenejte výkon, který máte nastavený pro daný setup v průměrném použití; je možné takto zaznamenat i příliš vysoké hodnoty nebo naopak velmi nízké. Pokud vyjde hodnocení vapu dobře, aplikace bude schopna určit optimální výkon pro ten či onen build.'
$ws.Range("A446").Value = 'cs'
$ws.Range("B446").Value = 'lab.vape.tc.label'
$ws.Range("C446").Value = 'Teplota'
$ws.Range("A447").Value = 'cs'
$ws.Range("B447").Value = 'lab.vape.tc.label.tooltip'
$ws.Range("C447").Value = 'Pokud máte spirálku, která podporuje režim teploty, je možné si zde zaznamenat vliv teploty na chuť a vůbec zážitek z vapování, včetně těch nepříjemných.'
$ws.Range("A448").Value = 'cs'
$ws.Range("B448").Value = 'lab.vape.airflow.label'
$ws.Range("C448").Value = 'Airflow'
$ws.Range("A449").Value = 'cs'
$ws.Range("B449").Value = 'lab.vape.airflow.label.tooltip'
$ws.Range("C449").Value = 'Jedná se o nastavení airflow na atomizéru (pokud jej umí); nula je úplně zavřeno nebo velmi tuhé MTL, pětka je plně otevřeno na velmi volné DL.'
$ws.Range("A450").Value = 'cs'
$ws.Range("B450").Value = 'lab.vape.juice.label'
$ws.Range("C450").Value = 'Juice flow'
$ws.Range("A451").Value = 'cs'
$ws.Range("B451").Value = 'lab.vape.juice.label.tooltip'
$ws.Range("C451").Value = 'Pokud atomizér podporuje juice flow, je možné si zde poznamenat optimální nastavení, aby nedocházelo k únikům. Nula je takřka zavřeno, pět je plně otevřeno.'
$ws.Range("A452").Value = 'cs'
$ws.Range("B452").Value = 'lab.vape.vape.title'
$ws.Range("C452").Value = 'DL/MTL'
$ws.Range("A453").Value = 'cs'
$ws.Range("B453").Value = 'lab.vape.mtl.label'
$ws.Range("C453").Value = 'Hodnocení MTL'
$ws.Range("A454").Value = 'cs'
$ws.Range("B454").Value = 'lab.vape.mtl.label.tooltip'
$ws.Range("C454").Value = 'Zde se jedná o kombinaci hodnocení chuťového projevu v MTL; čím vyšší hodnocení, tím lépe build funguje v MTL.'
$ws.Range("A455").Value = 'cs'
$ws.Range("B455").Value = 'lab.vape.dl.label'
$ws.Range("C455").Value = 'Hodnocení DL'
$ws.Range("A456").Value = 'cs'
$ws.Range("B456").Value = 'lab.vape.dl.label.tooltip'
$ws.Range("C456").Value = 'Hodnocení buildu pro DL; čím vyšší hodnocení, tím lépe build funguje v DL.'
$ws.Range("A457").Value = 'cs'
$ws.Range("B457").Value = 'lab.vape.clouds.label'
$ws.Range("C457").Value = 'Oblaka'
$ws.Range("A458").Value = 'cs'
$ws.Range("B458").Value = 'lab.vape.clouds.label.tooltip'
$ws.Range("C458").Value = 'Prosté hodnocení generovaných oblak; nízké hodnocení je komorní MTL mezi lidi, maximální hodnocení je prasostroj někde venku. Nebo mezi vapery.'
$ws.Range("A459").Value = 'cs'
$ws.Range("B459").Value = 'lab.vape.leaks.label'
$ws.Range("C459").Value = 'Úniky'
$ws.Range("A460").Value = 'cs'
$ws.Range("B460").Value = 'lab.vape.leaks.label.tooltip'
$ws.Range("C460").Value = 'Touto hodnotou je míněno, jakou tendenci má liquid unikat; toto je obecně nežádoucí - čím vyšší číslo, tím více má tank tendenci protékat. Nula naopak znamená žádné protečení.'
$ws.Range("A461").Value = 'cs'
$ws.Range("B461").Value = 'lab.vape.dryhit.label'
$ws.Range("C461").Value = 'Dryhit'
$ws.Range("A462").Value = 'cs'
$ws.Range("B462").Value = 'lab.vape.dryhit.label.tooltip'
$ws.Range("C462").Value = 'Touto hodnotou je míněno, jak moc je možné atomizér trápit, než se dostaví dryhit; vyšší hodnota obecně znamená stabilnější dodávání liquidu a bezproblémové bafání za sebou, nižší naopak vyžaduje střídmější a hodnoty k nule znamenají, že je něco špatně.'

# Match the "import" cell style (wrap text + 10pt font) used throughout the sheet
$ws.Range("A415:C462").WrapText = $true
$ws.Range("A415:C462").Font.Size = 10

# Rows whose wrapped text needs more than one line get an explicit row height
$ws.Rows.Item(428).RowHeight = 26.25
$ws.Rows.Item(430).RowHeight = 26.25
$ws.Rows.Item(433).RowHeight = 26.25
$ws.Rows.Item(439).RowHeight = 39
$ws.Rows.Item(441).RowHeight = 39
$ws.Rows.Item(445).RowHeight = 51.75
$ws.Rows.Item(447).RowHeight = 26.25
$ws.Rows.Item(449).RowHeight = 26.25
$ws.Rows.Item(451).RowHeight = 26.25
$ws.Rows.Item(454).RowHeight = 26.25
$ws.Rows.Item(458).RowHeight = 26.25
$ws.Rows.Item(460).RowHeight = 39
$ws.Rows.Item(462).RowHeight = 39

# Restore the selection/scroll state recorded by Excel after the edit
[void]$ws.Range("B460").Select()
